{"js": "// Replace the division problems in the first table with their updated\n// values. The table has 5 \"data\" rows (0, 4, 8, 12, 16) each holding 5\n// cells with text like \"86\u00f76=\" and 15 blank spacer rows in between.\n// We address every cell positionally (row, col) so duplicate values\n// (e.g. \"17\u00f72=\" appears twice, each replaced with a different number)\n// are each updated correctly instead of via a blind find/replace.\n\nconst table = context.document.body.tables.getFirst();\n\nconst updates = [\n  // [row, col, newValue]\n  [0, 0, \"86\u00f76=\"],\n  [0, 1, \"66\u00f76=\"],\n  [0, 2, \"89\u00f75=\"],\n  [0, 3, \"29\u00f73=\"],\n  [0, 4, \"55\u00f72=\"],\n\n  [4, 0, \"26\u00f73=\"],\n  [4, 1, \"85\u00f77=\"],\n  [4, 2, \"91\u00f72=\"],\n  [4, 3, \"84\u00f75=\"],\n  [4, 4, \"90\u00f75=\"],\n\n  [8, 0, \"33\u00f74=\"],\n  [8, 1, \"80\u00f74=\"],\n  [8, 2, \"90\u00f73=\"],\n  [8, 3, \"40\u00f72=\"],\n  [8, 4, \"52\u00f77=\"],\n\n  [12, 0, \"23\u00f74=\"],\n  [12, 1, \"10\u00f76=\"],\n  [12, 2, \"99\u00f77=\"],\n  [12, 3, \"24\u00f78=\"],\n  [12, 4, \"39\u00f75=\"],\n\n  [16, 0, \"97\u00f74=\"],\n  [16, 1, \"32\u00f75=\"],\n  [16, 2, \"95\u00f73=\"],\n  [16, 3, \"79\u00f74=\"],\n  [16, 4, \"66\u00f74=\"],\n];\n\nfor (const [row, col, value] of updates) {\n  table.getCell(row, col).value = value;\n}\n\nawait context.sync();\n", "ps1": "# Replace the division problems in the first table with their updated\n# values. The table has 5 \"data\" rows (Word 1-based rows 1, 5, 9, 13, 17)\n# each holding 5 cells with text like \"86\u00f76=\" and 15 blank spacer rows in\n# between. Cells are addressed positionally (row, col) so duplicate values\n# (e.g. \"17\u00f72=\" appears twice, each replaced with a different number) are\n# each updated correctly instead of via a blind find/replace.\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$updates = @(\n    @(1, 1, \"86\u00f76=\"),\n    @(1, 2, \"66\u00f76=\"),\n    @(1, 3, \"89\u00f75=\"),\n    @(1, 4, \"29\u00f73=\"),\n    @(1, 5, \"55\u00f72=\"),\n\n    @(5, 1, \"26\u00f73=\"),\n    @(5, 2, \"85\u00f77=\"),\n    @(5, 3, \"91\u00f72=\"),\n    @(5, 4, \"84\u00f75=\"),\n    @(5, 5, \"90\u00f75=\"),\n\n    @(9, 1, \"33\u00f74=\"),\n    @(9, 2, \"80\u00f74=\"),\n    @(9, 3, \"90\u00f73=\"),\n    @(9, 4, \"40\u00f72=\"),\n    @(9, 5, \"52\u00f77=\"),\n\n    @(13, 1, \"23\u00f74=\"),\n    @(13, 2, \"10\u00f76=\"),\n    @(13, 3, \"99\u00f77=\"),\n    @(13, 4, \"24\u00f78=\"),\n    @(13, 5, \"39\u00f75=\"),\n\n    @(17, 1, \"97\u00f74=\"),\n    @(17, 2, \"32\u00f75=\"),\n    @(17, 3, \"95\u00f73=\"),\n    @(17, 4, \"79\u00f74=\"),\n    @(17, 5, \"66\u00f74=\")\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $val = $u[2]\n    $t.Cell($row, $col).Range.Text = $val\n}\n"}
